# This workbook contains a Power-Query-backed table ("쿼리2") pulling live
# BJ-ranking data. The author re-ran "Refresh All" in Excel, which pulled
# fresh numbers from the source and stamped a new refresh time on every row.
# Since the external source isn't reachable here, we reproduce the refreshed
# worksheet state directly: updated point totals for the rows whose rank
# changed, a new refresh timestamp for every data row, the cursor left where
# the user clicked afterwards, and the auto-fit width nudge on the refresh-time
# column that Excel applies once the new values are in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refreshed "월별 누적별풍선" (monthly cumulative balloon) totals ---
# Only the rows whose underlying value actually moved between refreshes.
$ws.Range("C4").Value = 542371
$ws.Range("C5").Value = 536532
$ws.Range("C8").Value = 311114
$ws.Range("C10").Value = 203517

# --- New refresh timestamp, written to every data row (D2:D12) ---
$ws.Range("D2:D12").Value = 46018.011015520831

# --- Column D ("새로고침시간") re-autofits slightly narrower after refresh ---
$ws.Columns.Item(4).ColumnWidth = 17

# --- Cursor position left where the user clicked after refreshing ---
$ws.Range("F14").Select()
